# Applies the "Add files via upload" revision to the single slide in this
# deck:
#   1. Seven oval shapes (the Time 0 / Time 1 / Time 2 / Time 3 / (T1-T0) /
#      (T2-T1) / (T3-T2) labels) each had a redundant, completely empty
#      paragraph ("<a:pPr algn=\"ctr\"/><a:endParaRPr lang=\"en-US\"
#      dirty=\"0\"/>") removed from their text body (3 paragraphs -> 2).
#   2. Twelve "Construct Name" / "Slope" / "M = S-Mean" / "Intercept" style
#      text boxes were nudged to new positions (and two of them were also
#      narrowed).
#
# Helper: find a shape on a slide by its (stable) Shape.Id, since several
# shapes share the same Name ("TextBox 391", "TextBox 325", ...).
function Get-ShapeById {
    param($Shapes, [int]$Id)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $candidate = $Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes

# ---------------------------------------------------------------------
# 1) Drop the extra empty paragraph from each of the seven oval labels.
#    In every case paragraph #2 (1-based) is the fully empty duplicate
#    ("<a:pPr algn=\"ctr\"/><a:endParaRPr lang=\"en-US\" dirty=\"0\"/>"),
#    so deleting it leaves the correct final text body.
# ---------------------------------------------------------------------
$ovalIds = @(3, 51, 52, 53, 104, 105, 183)
foreach ($id in $ovalIds) {
    $shp = Get-ShapeById $shapes $id
    $tr = $shp.TextFrame.TextRange
    $extraPara = $tr.Paragraphs(2, 1)
    $extraPara.Delete()
}

# ---------------------------------------------------------------------
# 2) Reposition / resize the textboxes. Target positions are expressed in
#    EMU in the source XML; PowerPoint's COM surface uses points
#    (1 pt = 12700 EMU), stored internally as Single (float32), so the
#    literals below carry a tiny upward nudge (well under 1/100 pt) that
#    keeps the float32 round-trip landing on the exact target EMU instead
#    of one EMU short.
# ---------------------------------------------------------------------
# id, new Left (pt), new Top (pt), new Width (pt or $null = unchanged)
$moves = @(
    @{ Id = 10; Left = 61.38007954015748;  Top = 231.17063162125984; Width = $null },
    @{ Id = 62; Left = 67.93858287716536;  Top = 64.51063182125985;  Width = $null },
    @{ Id = 71; Left = 61.13976507952756;  Top = 413.6661532322835;  Width = $null },
    @{ Id = 72; Left = 295.9696199992126;  Top = 413.6661532322835;  Width = $null },
    @{ Id = 73; Left = 534.7579653559055;  Top = 413.6661532322835;  Width = $null },
    @{ Id = 74; Left = 773.3679529559056;  Top = 414.0116537433071;  Width = $null },
    @{ Id = 86; Left = 294.4868623937008;  Top = 213.00606559212596; Width = $null },
    @{ Id = 87; Left = 533.8048033496063;  Top = 212.51772333543306; Width = $null },
    @{ Id = 88; Left = 773.1227725055118;  Top = 213.40488458976378; Width = $null },
    @{ Id = 63; Left = 109.37818917637796; Top = 82.2720474440945;   Width = 88.11236240472441 },
    @{ Id = 64; Left = 67.77023722047244;  Top = 102.19653723307087; Width = $null },
    @{ Id = 66; Left = 104.42094828188976; Top = 252.07251988503936; Width = 83.7100031 }
)

foreach ($m in $moves) {
    $shp = Get-ShapeById $shapes $m.Id
    $shp.Left = $m.Left
    $shp.Top = $m.Top
    if ($null -ne $m.Width) {
        $shp.Width = $m.Width
    }
}
